$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number. Every data row (2-426)
# was bumped from 46074 (2026-02-21) to 46075 (2026-02-22).
$ws.Range("C2:C426").Value = 46075
